$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: the live-SLR data file moved from NonOnco_SLR_Data.xlsx to NonOnco_SLR_Import_Data.xlsx
$ws.Cells.Item(9, 3).Value = "\Testdata\Non_Oncology\DataFiles\LiveSLRPage\NonOnco_SLR_Import_Data.xlsx"

# New row 12: live-SLR report comparison data
$ws.Cells.Item(12, 1).Value = "test"
$ws.Cells.Item(12, 3).Value = "\Testdata\Non_Oncology\DataFiles\LiveSLRPage\NonOnco_SLRReport_Data.xlsx"
$ws.Cells.Item(12, 2).Value = "nononcology_liveslr_report_data"

[void]$ws.Range("B12").Select()
